$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.663.29"
Set-TextValue $ws.Range("E2") "  -0.68%  "
Set-TextValue $ws.Range("D3") "1.740.27"
Set-TextValue $ws.Range("E3") "  -1.63%  "
Set-TextValue $ws.Range("D4") "1.009"
Set-TextValue $ws.Range("E4") "  +1.07%  "
Set-TextValue $ws.Range("D5") "331.54"
Set-TextValue $ws.Range("E5") "  -1.01%  "
Set-TextValue $ws.Range("D6") "1.007"
Set-TextValue $ws.Range("E6") "  +1.06%  "
Set-TextValue $ws.Range("D7") "0.3827"
Set-TextValue $ws.Range("E7") "  +0.10%  "
Set-TextValue $ws.Range("D8") "0.3344"
Set-TextValue $ws.Range("E8") "  -2.46%  "
Set-TextValue $ws.Range("D9") "45.61"
Set-TextValue $ws.Range("E9") "  -2.47%  "
Set-TextValue $ws.Range("D10") "1.095"
Set-TextValue $ws.Range("E10") "  -3.99%  "
Set-TextValue $ws.Range("D11") "0.07133"
Set-TextValue $ws.Range("E11") "  -3.78%  "
Set-TextValue $ws.Range("D12") "1.008"
Set-TextValue $ws.Range("E12") "  +1.23%  "
Set-TextValue $ws.Range("D13") "21.99"
Set-TextValue $ws.Range("E13") "  -1.93%  "
Set-TextValue $ws.Range("D14") "6.094"
Set-TextValue $ws.Range("E14") "  -4.10%  "
Set-TextValue $ws.Range("D15") "1.746.30"
Set-TextValue $ws.Range("E15") "  -1.11%  "
Set-TextValue $ws.Range("D16") "6.987"
Set-TextValue $ws.Range("E16") "  -1.51%  "
Set-TextValue $ws.Range("D17") "0.00001047"
Set-TextValue $ws.Range("E17") "  -2.70%  "
Set-TextValue $ws.Range("D18") "0.06586"
Set-TextValue $ws.Range("E18") "  -1.17%  "
Set-TextValue $ws.Range("B19") "Dai"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D19") "1.005"
Set-TextValue $ws.Range("E19") "  +0.73%  "
Set-TextValue $ws.Range("B20") "Litecoin"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D20") "78.37"
Set-TextValue $ws.Range("E20") "  -4.79%  "
Set-TextValue $ws.Range("D21") "16.58"
Set-TextValue $ws.Range("E21") "  -4.64%  "
Set-TextValue $ws.Range("D22") "6.127"
Set-TextValue $ws.Range("E22") "  -4.58%  "
Set-TextValue $ws.Range("D23") "27.658.51"
Set-TextValue $ws.Range("E23") "  -0.65%  "
Set-TextValue $ws.Range("D24") "11.41"
Set-TextValue $ws.Range("E24") "  -5.47%  "
Set-TextValue $ws.Range("D25") "2.415"
Set-TextValue $ws.Range("E25") "  +1.56%  "
Set-TextValue $ws.Range("D26") "152.84"
Set-TextValue $ws.Range("E26") "  -0.16%  "
Set-TextValue $ws.Range("D27") "19.54"
Set-TextValue $ws.Range("E27") "  -5.66%  "
Set-TextValue $ws.Range("D28") "2.254"
Set-TextValue $ws.Range("E28") "  -6.93%  "
Set-TextValue $ws.Range("D29") "1.944.83"
Set-TextValue $ws.Range("E29") "  -1.11%  "
Set-TextValue $ws.Range("D30") "1.253"
Set-TextValue $ws.Range("E30") "  -13.10%  "
Set-TextValue $ws.Range("D31") "129.47"
Set-TextValue $ws.Range("E31") "  -3.50%  "
Set-TextValue $ws.Range("D32") "4.016"
Set-TextValue $ws.Range("E32") "  +1.59%  "
Set-TextValue $ws.Range("D33") "5.712"
Set-TextValue $ws.Range("E33") "  -6.90%  "
Set-TextValue $ws.Range("D34") "0.08644"
Set-TextValue $ws.Range("E34") "  -1.64%  "
Set-TextValue $ws.Range("D35") "11.87"
Set-TextValue $ws.Range("E35") "  -7.02%  "
Set-TextValue $ws.Range("D36") "1.533"
Set-TextValue $ws.Range("E36") "  +1.31%  "
Set-TextValue $ws.Range("B37") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D37") "5.070"
Set-TextValue $ws.Range("E37") "  -5.01%  "
Set-TextValue $ws.Range("B38") "TheSandbox"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D38") "0.6410"
Set-TextValue $ws.Range("E38") "  -5.98%  "
Set-TextValue $ws.Range("D39") "0.02243"
Set-TextValue $ws.Range("E39") "  -7.90%  "
Set-TextValue $ws.Range("D40") "0.06032"
Set-TextValue $ws.Range("E40") "  -4.78%  "
Set-TextValue $ws.Range("D41") "0.2075"
Set-TextValue $ws.Range("E41") "  -5.10%  "
Set-TextValue $ws.Range("D42") "1.191"
Set-TextValue $ws.Range("E42") "  -3.78%  "
Set-TextValue $ws.Range("B43") "Frax"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D43") "1.006"
Set-TextValue $ws.Range("E43") "  +0.96%  "
Set-TextValue $ws.Range("B44") "FraxShare"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "7.919"
Set-TextValue $ws.Range("E44") "  -4.21%  "
Set-TextValue $ws.Range("D45") "13.56"
Set-TextValue $ws.Range("E45") "  -4.01%  "
Set-TextValue $ws.Range("D46") "3.804"
Set-TextValue $ws.Range("E46") "  -0.85%  "
Set-TextValue $ws.Range("D47") "0.5926"
Set-TextValue $ws.Range("E47") "  -5.57%  "
Set-TextValue $ws.Range("D48") "126.00"
Set-TextValue $ws.Range("E48") "  -4.38%  "
Set-TextValue $ws.Range("D49") "1.966"
Set-TextValue $ws.Range("E49") "  -5.67%  "
Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.06930"
Set-TextValue $ws.Range("E50") "  -6.04%  "
Set-TextValue $ws.Range("B51") "EOS"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D51") "1.141"
Set-TextValue $ws.Range("E51") "  +0.01%  "
